$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.097.52"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.78"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.65"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.68"
$ws.Range("E7").Value = "  +9.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0796"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.15"
$ws.Range("E12").Value = "  +7.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.351.86"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.838"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("E15").Value = "  +10.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.055.31"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.47"
$ws.Range("E17").Value = "  +29.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.052.35"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.98"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.41"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.44"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +11.22%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.48"
$ws.Range("E26").Value = "  +4.59%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.64"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +8.76%  "
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.58"
$ws.Range("E33").Value = "  +5.05%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.110"
$ws.Range("E38").Value = "  +4.71%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  +15.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.32"
$ws.Range("E41").Value = "  +21.58%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.68"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.58"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.297.16"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  -19.66%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.89"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.239.25"
$ws.Range("E51").Value = "  -0.09%  "
